# Apply the cryptos list price/volume update (GitHub Actions commit, Mon Nov 13 14:25:42 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.916.62'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '2.063.51'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''245.96'
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("D6").Value = '''0.655'
$ws.Range("E6").Value = '  -1.84%  '
$ws.Range("D7").Value = '''58.93'
$ws.Range("E7").Value = '  -1.16%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''59.16'
$ws.Range("E9").Value = '  -2.05%  '
$ws.Range("D10").Value = '''0.372'
$ws.Range("E10").Value = '  -3.71%  '
$ws.Range("D11").Value = '''0.0780'
$ws.Range("E11").Value = '  -1.21%  '
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").Value = '''15.21'
$ws.Range("E13").Value = '  -4.21%  '
$ws.Range("D14").Value = '''0.880'
$ws.Range("E14").Value = '  +4.70%  '
$ws.Range("D15").Value = '2.359.61'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '''5.59'
$ws.Range("E16").Value = '  -3.39%  '
$ws.Range("D17").Value = '2.094.07'
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("D18").Value = '36.910.43'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").Value = '''17.55'
$ws.Range("E19").Value = '  -3.43%  '
$ws.Range("D20").Value = '''73.52'
$ws.Range("E20").Value = '  -2.54%  '
$ws.Range("D21").Value = '0.0₃0893'
$ws.Range("E21").Value = '  -1.17%  '
$ws.Range("D22").Value = '''5.44'
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("D23").Value = '''236.36'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("D26").Value = '''9.99'
$ws.Range("E26").Value = '  +5.71%  '
$ws.Range("D27").Value = '''2.21'
$ws.Range("E27").Value = '  +1.54%  '
$ws.Range("D28").Value = '''168.67'
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = '''5.47'
$ws.Range("E30").Value = '  +13.20%  '
$ws.Range("E31").Value = '  -0.99%  '
$ws.Range("E32").Value = '  +1.89%  '
$ws.Range("D33").Value = '''4.85'
$ws.Range("E33").Value = '  +6.38%  '
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").Value = '''2.36'
$ws.Range("E35").Value = '  +3.14%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  +5.37%  '
$ws.Range("D38").Value = '''0.0854'
$ws.Range("E38").Value = '  -5.89%  '
$ws.Range("E39").Value = '  -2.35%  '
$ws.Range("D40").Value = '''0.0222'
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").Value = '''2.97'
$ws.Range("E41").Value = '  -6.34%  '
$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").Value = '''4.92'
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '''1.16'
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = '''0.0964'
$ws.Range("E44").Value = '  -9.77%  '
$ws.Range("D45").Value = '''97.04'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("E46").Value = '  -4.36%  '
$ws.Range("D47").Value = '1.317.85'
$ws.Range("E47").Value = '  +2.08%  '
$ws.Range("E48").Value = '  -3.44%  '
$ws.Range("D49").Value = '''2.87'
$ws.Range("E49").Value = '  -1.97%  '
$ws.Range("D50").Value = '''6.86'
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("D51").Value = '2.247.00'
$ws.Range("E51").Value = '  +0.27%  '
